# Commit: "docs: Remove references to Visual Recognition API keys"
#
# The "Set up a class account" guide used to contain two list items about
# creating/adding Watson Visual Recognition API credentials. Those two
# whole list items (paragraphs) are removed; everything else in the list
# (Set up accounts for your students / Prepare a lesson plan / Try the
# worksheets out for yourself / Check your group settings / If you run
# into any problems...) is left as-is, just shifted up.
#
# The "Last updated" date stamped in the document footer is also bumped.

$d = $word.ActiveDocument

# --- 1. Remove the "Create ... Visual Recognition credentials ..." and
#        "Add ... Visual Recognition credentials ..." list items -------

$startFind = $d.Content
$foundStart = $startFind.Find.Execute(
    "Create Watson Visual Recognition credentials for your group to use",
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart) {
    $startPos = $startFind.Start

    $endFind = $d.Range($startPos, $d.Content.End)
    $foundEnd = $endFind.Find.Execute(
        "Set up accounts for your students",
        $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

    if ($foundEnd) {
        $endPos = $endFind.Start
        $deleteRange = $d.Range($startPos, $endPos)
        $deleteRange.Delete()
    }
}

# --- 2. Bump the "Last updated" date in the footer ---------------------

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)
    $footer = $section.Footers(1)
    if ($footer.Exists) {
        $footer.Range.Find.Execute(
            "1 October 2020", $false, $false, $false, $false, $false,
            $true, 1, $false, "9 February 2021", 2)
    }
}
